# The source workbook has a single sheet "2025-02-17" whose A1:B1 cells
# carry a bold/bordered header style (s=1) but no values.
#
# The target revision:
#   - adds a new sheet "2025-02-18" right after "2025-02-17"
#   - gives that new sheet's A1/B1 the same header style, with the text
#     "Цена за м2" (A1) and "Координаты" (B1) -- these become new shared
#     strings
#   - empties the original "2025-02-17" sheet's A1:B1 back to blank cells
#   - leaves "2025-02-17" as the active/selected sheet

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new sheet right after the existing one.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "2025-02-18"

# Carry over the existing header formatting (bold, centered, bordered)
# from sheet 1's A1:B1 onto the new sheet's A1:B1 before filling it in.
$ws1.Range("A1:B1").Copy()
$ws2.Range("A1:B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws2.Range("A1").Value = "Цена за м2"
$ws2.Range("B1").Value = "Координаты"

# The original sheet loses its (previously empty/styled) header cells.
$ws1.Range("A1:B1").Clear()

# Keep the first sheet the active one, as in the source workbook.
$ws1.Activate()
